$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Ccl12 -> Ccr3 -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ccr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 3.288126333333333
$ws.Range("H2").Value = 9.864379
$ws.Range("I2").Value = 0.05813306630866938
$ws.Range("J2").Value = 0.05813306630866937
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.438062
$ws.Range("N2").Value = 1.314186
$ws.Range("O2").Value = 0.6074000808827777
$ws.Range("P2").Value = 0.6074000808827777
$ws.Range("Q2").Value = 1.440403197832667
$ws.Range("R2").Value = 12.963628780494
$ws.Range("S2").Value = 0.03531002917784967
$ws.Range("T2").Value = 0.03531002917784966

# Row 3: ECs -> Ccl12 -> Ccr3 -> M2
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ccr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 3.288126333333333
$ws.Range("H3").Value = 9.864379
$ws.Range("I3").Value = 0.05813306630866938
$ws.Range("J3").Value = 0.05813306630866937
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2831463333333333
$ws.Range("N3").Value = 0.8494390000000001
$ws.Range("O3").Value = 0.3925999191172223
$ws.Range("P3").Value = 0.3925999191172223
$ws.Range("Q3").Value = 0.9310209148201111
$ws.Range("R3").Value = 8.379188233381001
$ws.Range("S3").Value = 0.02282303713081972
$ws.Range("T3").Value = 0.02282303713081971

# Row 4: M2 -> Ccl12 -> Ccr3 -> ECs
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ccr3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 53.27393966666666
$ws.Range("H4").Value = 159.821819
$ws.Range("I4").Value = 0.9418669336913307
$ws.Range("J4").Value = 0.9418669336913306
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.438062
$ws.Range("N4").Value = 1.314186
$ws.Range("O4").Value = 0.6074000808827777
$ws.Range("P4").Value = 0.6074000808827777
$ws.Range("Q4").Value = 23.33728855825933
$ws.Range("R4").Value = 210.035597024334
$ws.Range("S4").Value = 0.5720900517049281
$ws.Range("T4").Value = 0.5720900517049281

# Row 5: M2 -> Ccl12 -> Ccr3 -> M2
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Ccl12"
$ws.Range("C5").Value = "Ccr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 53.27393966666666
$ws.Range("H5").Value = 159.821819
$ws.Range("I5").Value = 0.9418669336913307
$ws.Range("J5").Value = 0.9418669336913306
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2831463333333333
$ws.Range("N5").Value = 0.8494390000000001
$ws.Range("O5").Value = 0.3925999191172223
$ws.Range("P5").Value = 0.3925999191172223
$ws.Range("Q5").Value = 15.08432067883789
$ws.Range("R5").Value = 135.758886109541
$ws.Range("S5").Value = 0.3697768819864026
$ws.Range("T5").Value = 0.3697768819864025
